$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# --- Row 2 (the "22/10/21" meeting row) ---

# 1. Date cell: "22/10/21" -> "28/10/21"
$cDate = $t.Cell(2, 1)
$cDate.Range.Find.Execute("22/10/21", $false, $false, $false, $false, $false, $true, 1, $false, "28/10/21", 2) | Out-Null

# 2. Agenda cell: "Group assignment meeting 2" -> "Group assignment meeting 3"
$cAgenda = $t.Cell(2, 2)
$cAgenda.Range.Find.Execute("Group assignment meeting 2", $false, $false, $false, $false, $false, $true, 1, $false, "Group assignment meeting 3", 2) | Out-Null

# 3. Outcome cell: populate the empty bulleted paragraph with 5 bullet points
$cOutcome = $t.Cell(2, 3)
$pOutcome = $cOutcome.Range.Paragraphs.Item(1)
$bullets = "Looked at each others use case diagrams and discussed the resolving of conflicts between them.`r" + `
    "Discussed the use of <<include>> and <<extend>>`r" + `
    "Begun work on collating all use case diagrams into one that includes actors from all our subsystems`r" + `
    "Clarified with each other how  to do the activity diagrams and went over examples in the slides`r" + `
    "Agenda includes finishing activity, class and sequence diagrams for next week so that we can begin task 2."
$pOutcome.Range.Text = $bullets

# 4. Time started cell: "4:00 pm" -> "3:00pm"
$cStart = $t.Cell(2, 4)
$cStart.Range.Find.Execute("4:00 pm", $false, $false, $false, $false, $false, $true, 1, $false, "3:00pm", 2) | Out-Null

# 5. Time ended cell: "-" -> "4:00pm"
$cEnd = $t.Cell(2, 5)
$cEnd.Range.Find.Execute("-", $false, $false, $false, $false, $false, $true, 1, $false, "4:00pm", 2) | Out-Null

